$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Enter the attendance data ("T1") in column B of the "Ninja" sheet.
#    This feeds the existing P-column formula
#    (=IFERROR(SUM(B:O)/COUNT(B:O)*100,0)) which recalculates to 100/0.
# ---------------------------------------------------------------------------
$wsNinja = $wb.Worksheets.Item("Ninja")

$attendance = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 1
}

foreach ($row in $attendance.Keys | Sort-Object) {
    $wsNinja.Cells.Item($row, 2).Value = $attendance[$row]
}

# ---------------------------------------------------------------------------
# 2. Reset the scroll position ("topLeftCell") of the "Astronauta" sheet back
#    to A1, leaving its current selection untouched, without disturbing which
#    sheet tab is active.
# ---------------------------------------------------------------------------
$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 3. Re-activate "Ninja" (the originally active tab), reset its scroll
#    position, zoom in from 70% to 85%, and move the selection from
#    P27:P28 to A16.
# ---------------------------------------------------------------------------
$wsNinja.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 85
$wsNinja.Range("A16").Select()
